$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (right after the header), shifting existing data down.
$ws.Rows.Item(2).Insert()

# Fill the new row with the latest date and the same price values as the rest of the data.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-02-10"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
